$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) - force text format so Excel doesn't
# auto-convert month/year strings into date serial numbers
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "April 2025"
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "May 2025"

# Update data values (row 2)
$ws.Range("A2").Value = 2.078
$ws.Range("B2").Value = -0.098
$ws.Range("C2").Value = 0.122
$ws.Range("D2").Value = 0.132
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = -0.808
$ws.Range("G2").Value = 1.421
